$d = $word.ActiveDocument

function Replace-Once($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output ("WARNING: replace failed for: " + $findText)
    }
    return $ok
}

# --- Air Quality paragraph: collapse an accidental double space around the
#     two ZZairqualityXXXZZ merge placeholders ---
Replace-Once "to be `$ZZairqualityAvgZZ  per acre for a total value of" "to be `$ZZairqualityAvgZZ per acre for a total value of"
Replace-Once "`$ZZairqualityTotalZZ  annually" "`$ZZairqualityTotalZZ annually"

# --- Biodiversity paragraph: same double-space cleanup ---
Replace-Once "is `$ZZbiodiversityAvgZZ  per acre per year for a total value of" "is `$ZZbiodiversityAvgZZ per acre per year for a total value of"
Replace-Once "`$ZZbiodiversityTotalZZ  annually" "`$ZZbiodiversityTotalZZ annually"

# --- Carbon sequestration paragraph: same double-space cleanup ---
Replace-Once "to be `$ZZcarbonAvgZZ  per acre per year, for a total value of" "to be `$ZZcarbonAvgZZ per acre per year, for a total value of"

# --- Cultural value paragraph: same double-space cleanup ---
Replace-Once "`$ZZculturalTotalZZ  annually" "`$ZZculturalTotalZZ annually"

# --- Total economic value paragraph: double-space cleanup AND rename the
#     merge field from ZZtotalTotalZZ to ZZtotalTotalThousandZZ (only the
#     first occurrence, inside this sentence — not the one in the summary
#     table further down) ---
Replace-Once "to be `$ZZtotalAvgZZ  per acre, for a total value of approximately `$ZZtotalTotalZZ" "to be `$ZZtotalAvgZZ per acre, for a total value of approximately `$ZZtotalTotalThousandZZ"
Replace-Once "`$ZZtotalTotalThousandZZ  annually" "`$ZZtotalTotalThousandZZ annually"

# --- Footer: bump the template's last-updated date ---
foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("4/14/2020", $true, $false, $false, $false, $false, $true, 1, $false, "1/26/2021", 2) | Out-Null
        }
    }
}

$d.Content.Find.Execute("4/14/2020", $true, $false, $false, $false, $false, $true, 1, $false, "1/26/2021", 2) | Out-Null

Write-Output "done"
